# Adds a new "2022-Q3" sheet (right after the "总计" summary sheet) containing
# the fund-holding detail table for that quarter, and updates the "总计"
# summary sheet with a new leading row for 2022-Q3 (pushing the existing
# rows down by one).

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)

# Helper: write a value as *text* (matching the workbook's existing
# convention of storing numeric-looking figures as inline strings), without
# leaving a stray explicit number-format style behind on the cell.
function Set-TextCell($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Helper: write a plain (non-numeric-looking) text value, then stamp the
# bold/boxed header style copied from $styleCell onto it, without Excel
# synthesizing a brand-new merged style.
function Set-StyledTextCell($cell, [string]$text, $styleCell) {
    $cell.Value = $text
    $styleCell.Copy()
    $cell.PasteSpecial(-4122)
}

# Helper: write a number, then stamp the bold/boxed style copied from
# $styleCell onto it (used for the running-index column).
function Set-StyledNumberCell($cell, $num, $styleCell) {
    $cell.Value = $num
    $styleCell.Copy()
    $cell.PasteSpecial(-4122)
    $cell.Value = $num
}

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet right after "总计".
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

# Header row (style copied from the summary sheet's bold/boxed header cell).
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 2 + $i  # headers start at column B
    Set-StyledTextCell $q3.Cells.Item(1,$col) $headers[$i] $summary.Cells.Item(1,2)
}

# Data rows: code, name, scale, totalPosition, positionPct, marketValue, rank
$rows = @(
    @("002685", "中欧丰泓沪港深灵活配置混合A", "42.36", "92.77", "5.36", "2.2705", 8),
    @("002686", "中欧丰泓沪港深灵活配置混合C", "7.40",  "92.77", "5.36", "0.3966", 8),
    @("517880", "华泰柏瑞中证沪港深品牌消费50ETF", "0.47", "92.73", "3.45", "0.0162", 9)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = 2 + $i
    $data = $rows[$i]

    # Column A: sequential 0-based index, bold/boxed style like the summary sheet.
    Set-StyledNumberCell $q3.Cells.Item($r,1) $i $summary.Cells.Item(2,1)

    Set-TextCell $q3.Cells.Item($r,2) $data[0]
    Set-TextCell $q3.Cells.Item($r,3) $data[1]
    Set-TextCell $q3.Cells.Item($r,4) $data[2]
    Set-TextCell $q3.Cells.Item($r,5) $data[3]
    Set-TextCell $q3.Cells.Item($r,6) $data[4]
    Set-TextCell $q3.Cells.Item($r,7) $data[5]
    $q3.Cells.Item($r,8).Value = $data[6]
}

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a new row for 2022-Q3 and
#    renumber the running index column.
# ---------------------------------------------------------------------
$summary.Rows.Item(2).Insert()

$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 3
$summary.Cells.Item(2,4).Value = 2.68

# Clear the formatting the row-insert borrowed from the row below, then
# restore the proper boxed/bold style on column A (matching the other rows).
$summary.Cells.Item(2,2).Style = "Normal"
$summary.Cells.Item(2,3).Style = "Normal"
$summary.Cells.Item(2,4).Style = "Normal"

Set-StyledNumberCell $summary.Cells.Item(2,1) 0 $summary.Cells.Item(3,1)

# Renumber the existing rows' running index (they all shift down by one).
for ($r = 3; $r -le 9; $r++) {
    $summary.Cells.Item($r,1).Value = $r - 2
}
